$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.746.82'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '3.798.21'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.517'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.451'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000247'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("D14").Value = '4.439.08'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").Value = '3.827.72'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '67.772.75'
$ws.Range("E16").Value = '  +0.29%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.65%  '

$ws.Range("E23").Value = '  -4.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").Value = '3.948.73'
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.17%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0994'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.67%  '

$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.996'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.24%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '45.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.297'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.63'
$ws.Range("D46").Style = "Normal"

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.70%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '393.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.77%  '
